$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"22.75000000000012"
$ws.Range("G2").Value = [double]"1.159224076396903e-06"
$ws.Range("H2").Value = [double]"2.583161607423796e-05"
$ws.Range("I2").Value = [double]"1.110223024625157e-16"
$ws.Range("K2").Value = [double]"6.365825451238024"
$ws.Range("L2").Value = "[3.417922702382862, 9.313728200093186]"
$ws.Range("M2").Value = [double]"2.738413904257442e-05"
$ws.Range("N2").Value = [double]"5.476827808514884e-05"
$ws.Range("O2").Value = [double]"-1.232737057264464"
$ws.Range("P2").Value = "[-1.7484739893853103, -0.7170001251436169]"
$ws.Range("Q2").Value = [double]"3.645341768443444e-06"
$ws.Range("R2").Value = [double]"4.75336125305148e-06"
$ws.Range("S2").Value = [double]"13.98877356041467"
$ws.Range("T2").Value = "[12.321687018302462, 15.65586010252688]"
$ws.Range("W2").Value = [double]"4.463463463463487"
$ws.Range("X2").Value = [double]"2.596096096096111"
$ws.Range("Y2").Value = [double]"6.330830830830863"

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = [double]"24.18000000000034"
$ws.Range("G3").Value = [double]"2.845748369961765e-05"
$ws.Range("H3").Value = [double]"0.0001215996211675573"
$ws.Range("K3").Value = [double]"6.094906309119983"
$ws.Range("L3").Value = "[3.0553298273095084, 9.134482790930457]"
$ws.Range("M3").Value = [double]"9.980915084106812e-05"
$ws.Range("N3").Value = [double]"9.980915084106812e-05"
$ws.Range("O3").Value = [double]"1.490605523324887"
$ws.Range("P3").Value = "[0.8616580451287321, 2.1195530015210418]"
$ws.Range("Q3").Value = [double]"4.75336125305148e-06"
$ws.Range("R3").Value = [double]"4.75336125305148e-06"
$ws.Range("S3").Value = [double]"12.59400299706103"
$ws.Range("T3").Value = "[10.77317554840604, 14.414830445716017]"
$ws.Range("W3").Value = [double]"18.44360360360386"
$ws.Range("X3").Value = [double]"16.02318318318341"
$ws.Range("Y3").Value = [double]"20.86402402402432"
